# Slide 1's speaker notes ("Note from Chuck...") are translated/expanded
# into Greek. The notes-page body placeholder only supports plain-text
# assignment in this host, so we set the full resulting text in one shot.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$notesBody = $s.NotesPage.Shapes.Item(1)
$notesBody.TextFrame.TextRange.Text = "Σημείωση από τον  Chuck. Εάν χρησιμοποιείτε αυτό το υλικό, μπορείτε να αφαιρέσετε το λογότυπο UM και να το αντικαταστήσετε με το δικό σας, αλλά διατηρήστε το λογότυπο CC-BY στην πρώτη σελίδα καθώς την/τις σελίδα/ες αναγνώρισης."
